# Release QRPH.PCC 1.0.0 for TI
# - Bump Version metadata from "1.0.0-comment" to "1.0.0"
# - Update publication Date
# - Append "|4.0.1" FHIR version pins to every canonical()/Reference()/Quantity
#   type reference and ValueSet binding URL on the Elements sheet
# - Column widths on the Elements sheet grow to fit the longer Type(s) /
#   Binding Value Set text (best-fit re-autosize)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.0.0"
$meta.Range("B8").Value = "2025-10-02T10:26:47-05:00"

# ---------------------------------------------------------------------------
# Elements sheet - Type(s) column (K) and Binding Value Set column (Z)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

$ws.Range("Z6").Value = "http://hl7.org/fhir/ValueSet/languages|4.0.1"
$ws.Range("K12").Value = "canonical(PlanDefinition|4.0.1|Questionnaire|4.0.1|Measure|4.0.1|ActivityDefinition|4.0.1|OperationDefinition|4.0.1)`n"
$ws.Range("K14").Value = "Reference(CarePlan|4.0.1)`n"
$ws.Range("K15").Value = "Reference(CarePlan|4.0.1)`n"
$ws.Range("K16").Value = "Reference(CarePlan|4.0.1)`n"
$ws.Range("Z19").Value = "http://hl7.org/fhir/ValueSet/care-plan-category|4.0.1"
$ws.Range("K22").Value = "Reference(Patient|4.0.1|Group|4.0.1)`n"
$ws.Range("K23").Value = "Reference(Encounter|4.0.1)`n"
$ws.Range("K26").Value = "Reference(Patient|4.0.1|Practitioner|4.0.1|PractitionerRole|4.0.1|Device|4.0.1|RelatedPerson|4.0.1|Organization|4.0.1|CareTeam|4.0.1)`n"
$ws.Range("K27").Value = "Reference(Patient|4.0.1|Practitioner|4.0.1|PractitionerRole|4.0.1|Device|4.0.1|RelatedPerson|4.0.1|Organization|4.0.1|CareTeam|4.0.1)`n"
$ws.Range("K28").Value = "Reference(CareTeam|4.0.1)`n"
$ws.Range("K29").Value = "Reference(Condition|4.0.1)`n"
$ws.Range("K30").Value = "Reference(Resource|4.0.1)`n"
$ws.Range("K31").Value = "Reference(Goal|4.0.1)`n"
$ws.Range("Z36").Value = "http://hl7.org/fhir/ValueSet/care-plan-activity-outcome|4.0.1"
$ws.Range("K37").Value = "Reference(Resource|4.0.1)`n"
$ws.Range("K39").Value = "Reference(Appointment|4.0.1|CommunicationRequest|4.0.1|DeviceRequest|4.0.1|MedicationRequest|4.0.1|NutritionOrder|4.0.1|Task|4.0.1|ServiceRequest|4.0.1|VisionPrescription|4.0.1|RequestGroup|4.0.1)`n"
$ws.Range("K45").Value = "canonical(PlanDefinition|4.0.1|ActivityDefinition|4.0.1|Questionnaire|4.0.1|Measure|4.0.1|OperationDefinition|4.0.1)`n"
$ws.Range("Z48").Value = "http://hl7.org/fhir/ValueSet/procedure-code|4.0.1"
$ws.Range("Z49").Value = "http://hl7.org/fhir/ValueSet/clinical-findings|4.0.1"
$ws.Range("K50").Value = "Reference(Condition|4.0.1|Observation|4.0.1|DiagnosticReport|4.0.1|DocumentReference|4.0.1)`n"
$ws.Range("K51").Value = "Reference(Goal|4.0.1)`n"
$ws.Range("K56").Value = "Reference(Location|4.0.1)`n"
$ws.Range("K57").Value = "Reference(Practitioner|4.0.1|PractitionerRole|4.0.1|Organization|4.0.1|RelatedPerson|4.0.1|Patient|4.0.1|CareTeam|4.0.1|HealthcareService|4.0.1|Device|4.0.1)`n"
$ws.Range("K58").Value = "CodeableConcept`nReference(Medication|4.0.1|Substance|4.0.1)"
$ws.Range("Z58").Value = "http://hl7.org/fhir/ValueSet/medication-codes|4.0.1"
$ws.Range("K59").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$ws.Range("K60").Value = "Quantity {SimpleQuantity|4.0.1}`n"

# ---------------------------------------------------------------------------
# Re-autosize the two columns whose content grew (Type(s) = column 11,
# Binding Value Set = column 26) so their stored width matches the longer text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 172.34635416666666
$ws.Columns.Item(26).ColumnWidth = 49.936197916666664
